$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.76321629402796
$ws.Range("C2").Value = 6.977608749697993
$ws.Range("D2").Value = 10.27562151199105
$ws.Range("E2").Value = 10.16221892711302
$ws.Range("F2").Value = 53.19640514517793
$ws.Range("L2").Value = 10.07743624903725
$ws.Range("M2").Value = 17.6751976882505
$ws.Range("B3").Value = 20.56136534886759
$ws.Range("C3").Value = 6.545288562991816
$ws.Range("D3").Value = 10.15621449206209
$ws.Range("E3").Value = 10.13950400312873
$ws.Range("F3").Value = 52.081782991659
$ws.Range("L3").Value = 10.09209268637333
$ws.Range("M3").Value = 17.67393970846745
$ws.Range("B4").Value = 20.44576969684817
$ws.Range("C4").Value = 6.264247296778958
$ws.Range("D4").Value = 10.08222054588462
$ws.Range("E4").Value = 10.12499934053215
$ws.Range("F4").Value = 51.39215742805976
$ws.Range("L4").Value = 10.10261598828078
$ws.Range("M4").Value = 17.67828721713881
$ws.Range("B5").Value = 20.40081852552255
$ws.Range("C5").Value = 6.145800879673959
$ws.Range("D5").Value = 10.05191071914747
$ws.Range("E5").Value = 10.11894723115809
$ws.Range("F5").Value = 51.11011740943114
$ws.Range("L5").Value = 10.10728735383755
$ws.Range("M5").Value = 17.68134529523607
$ws.Range("B6").Value = 20.39348614657114
$ws.Range("C6").Value = 6.125896313450055
$ws.Range("D6").Value = 10.04686872945196
$ws.Range("E6").Value = 10.11793371346354
$ws.Range("F6").Value = 51.06323286136561
$ws.Range("L6").Value = 10.10808615715359
$ws.Range("M6").Value = 17.68193073778849
$ws.Range("B7").Value = 20.44515467217441
$ws.Range("C7").Value = 6.262665753540227
$ws.Range("D7").Value = 10.08181239341741
$ws.Range("E7").Value = 10.12491829344967
$ws.Range("F7").Value = 51.38835742474631
$ws.Range("L7").Value = 10.10267743738798
$ws.Range("M7").Value = 17.67832325266093
$ws.Range("B8").Value = 20.69192284605454
$ws.Range("C8").Value = 6.831781386851232
$ws.Range("D8").Value = 10.23459882119556
$ws.Range("E8").Value = 10.15450187114349
$ws.Range("F8").Value = 52.81335284625715
$ws.Range("L8").Value = 10.0821733121836
$ws.Range("M8").Value = 17.67370094145012
$ws.Range("B9").Value = 21.23924707961388
$ws.Range("C9").Value = 7.824177648357022
$ws.Range("D9").Value = 10.52823790623437
$ws.Range("E9").Value = 10.20814250218221
$ws.Range("F9").Value = 55.55317911498504
$ws.Range("L9").Value = 10.05406845544233
$ws.Range("M9").Value = 17.70527014309124
$ws.Range("B10").Value = 21.67600274356233
$ws.Range("C10").Value = 8.478250637465971
$ws.Range("D10").Value = 10.73953937687238
$ws.Range("E10").Value = 10.24495939471958
$ws.Range("F10").Value = 57.51628383688331
$ws.Range("L10").Value = 10.04081115358301
$ws.Range("M10").Value = 17.75319598954625
$ws.Range("B11").Value = 21.88129162367294
$ws.Range("C11").Value = 8.759620414251589
$ws.Range("D11").Value = 10.83453680030563
$ws.Range("E11").Value = 10.26116114932803
$ws.Range("F11").Value = 58.39528615126658
$ws.Range("L11").Value = 10.03638687428107
$ws.Range("M11").Value = 17.78034195012373
$ws.Range("B12").Value = 21.95989887169857
$ws.Range("C12").Value = 8.863852382217932
$ws.Range("D12").Value = 10.87033426203832
$ws.Range("E12").Value = 10.26721901117585
$ws.Range("F12").Value = 58.72587409507312
$ws.Range("L12").Value = 10.03494257431568
$ws.Range("M12").Value = 17.79138605564508
$ws.Range("B13").Value = 21.94293205826419
$ws.Range("C13").Value = 8.841507047008429
$ws.Range("D13").Value = 10.86263269005824
$ws.Range("E13").Value = 10.26591776526293
$ws.Range("F13").Value = 58.65478095806817
$ws.Range("L13").Value = 10.03524335246426
$ws.Range("M13").Value = 17.78897358077518
$ws.Range("B14").Value = 21.88774174614223
$ws.Range("C14").Value = 8.768242016561572
$ws.Range("D14").Value = 10.83748548706569
$ws.Range("E14").Value = 10.26166107746956
$ws.Range("F14").Value = 58.42253073841886
$ws.Range("L14").Value = 10.03626341953435
$ws.Range("M14").Value = 17.78123525927681
$ws.Range("B15").Value = 21.85404680120337
$ws.Range("C15").Value = 8.723063649627333
$ws.Range("D15").Value = 10.82205875991617
$ws.Range("E15").Value = 10.25904369005426
$ws.Range("F15").Value = 58.27996772153743
$ws.Range("L15").Value = 10.03691833422556
$ws.Range("M15").Value = 17.77659474316815
$ws.Range("B16").Value = 21.66271265883937
$ws.Range("C16").Value = 8.459537777294084
$ws.Range("D16").Value = 10.73330721397141
$ws.Range("E16").Value = 10.24388967743671
$ws.Range("F16").Value = 57.45853464763976
$ws.Range("L16").Value = 10.04113262462789
$ws.Range("M16").Value = 17.75152915820444
$ws.Range("B17").Value = 21.54697008781206
$ws.Range("C17").Value = 8.293737220753911
$ws.Range("D17").Value = 10.67856271682715
$ws.Range("E17").Value = 10.23445406729067
$ws.Range("F17").Value = 56.95083404100146
$ws.Range("L17").Value = 10.04412948726079
$ws.Range("M17").Value = 17.737518487566
$ws.Range("B18").Value = 21.48102482537235
$ws.Range("C18").Value = 8.196850407299589
$ws.Range("D18").Value = 10.64697010841545
$ws.Range("E18").Value = 10.22897544943037
$ws.Range("F18").Value = 56.65751231165109
$ws.Range("L18").Value = 10.04600441768186
$ws.Range("M18").Value = 17.72996324385313
$ws.Range("B19").Value = 21.4588070001769
$ws.Range("C19").Value = 8.163784419350574
$ws.Range("D19").Value = 10.63625579597304
$ws.Range("E19").Value = 10.22711160167348
$ws.Range("F19").Value = 56.55798236487698
$ws.Range("L19").Value = 10.04666520536776
$ws.Range("M19").Value = 17.72749172045719
$ws.Range("B20").Value = 21.55922679485278
$ws.Range("C20").Value = 8.311544508447311
$ws.Range("D20").Value = 10.68440133562119
$ws.Range("E20").Value = 10.2354638260526
$ws.Range("F20").Value = 57.00501667418031
$ws.Range("L20").Value = 10.04379481540629
$ws.Range("M20").Value = 17.73895788121139
$ws.Range("B21").Value = 21.90392955207732
$ws.Range("C21").Value = 8.789824530692481
$ws.Range("D21").Value = 10.84487671910196
$ws.Range("E21").Value = 10.26291346073618
$ws.Range("F21").Value = 58.49081184494358
$ws.Range("L21").Value = 10.03595752947919
$ws.Range("M21").Value = 17.78348747928734
$ws.Range("B22").Value = 22.1342385130235
$ws.Range("C22").Value = 9.088909640996826
$ws.Range("D22").Value = 10.9487243865137
$ws.Range("E22").Value = 10.28040284695567
$ws.Range("F22").Value = 59.44852022281044
$ws.Range("L22").Value = 10.0321823036225
$ws.Range("M22").Value = 17.81704386016261
$ws.Range("B23").Value = 22.01088549948896
$ws.Range("C23").Value = 8.930514261058068
$ws.Range("D23").Value = 10.8933980232923
$ws.Range("E23").Value = 10.27110926235287
$ws.Range("F23").Value = 58.93867418625543
$ws.Range("L23").Value = 10.03407396541892
$ws.Range("M23").Value = 17.79872822334587
$ws.Range("B24").Value = 21.5536836729179
$ws.Range("C24").Value = 8.303498711652555
$ws.Range("D24").Value = 10.68176206647118
$ws.Range("E24").Value = 10.23500748227889
$ws.Range("F24").Value = 56.98052516493243
$ws.Range("L24").Value = 10.04394564712904
$ws.Range("M24").Value = 17.7383055746701
$ws.Range("B25").Value = 21.08479358039419
$ws.Range("C25").Value = 7.568917841879596
$ws.Range("D25").Value = 10.4495225010029
$ws.Range("E25").Value = 10.19409191574789
$ws.Range("F25").Value = 54.81958752117841
$ws.Range("L25").Value = 10.06037428862328
$ws.Range("M25").Value = 17.69238019449641
